$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (old C -> E), shifting the
# "as-of" values column two places to the right and making room for two
# additional weekly snapshot columns.
$ws.Range("C1:D1").EntireColumn.Insert()

# Row 1 headers: D1 keeps the date that used to live in B1, B1 becomes
# the newest snapshot date, and C1 picks up the other newly-added
# snapshot date.
$ws.Range("D1").Value = "Jun_13"
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two new data columns (rows 2-27) with the same "UN" markers
# used in column B.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Keep the explicit width on column C (8.0 chars, same as the original
# "as-of" column) and carry it over to the two other "as-of" columns too.
$ws.Columns("C:E").ColumnWidth = 7.1666666
